# Clean UI: Remove System Status section and info messages for minimal interface
# (Data refresh on KRA_Database + Database_Summary sheets)

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("KRA_Database")
$wsSummary = $wb.Worksheets.Item("Database_Summary")

# --- KRA_Database sheet -----------------------------------------------
# Pre-format the preAmount column (E2:E7) as Text so the comma-formatted
# amount strings are not auto-coerced into numbers by Excel's smart entry.
$wsData.Range("E2:E7").NumberFormat = "@"

# Row 2
$wsData.Range("A2").Value = "29TH AUGUST, 2025"
$wsData.Range("B2").Value = "A001126762Z"
$wsData.Range("C2").Value = "Peter Kimutai Telengech"
$wsData.Range("D2").Value = "NOTICE UNDER SECTION 29 OF THE TAX PROCEDURES ACT, 2015"
$wsData.Range("E2").Value = "14,769.50"
$wsData.Range("H2").Value = "Franciscar Nyangweta"
$wsData.Range("J2").Value = "2025-09-26 11:25:36"
$wsData.Range("L2").Value = 6

# Row 3
$wsData.Range("A3").Value = "04th September, 2025"
$wsData.Range("B3").Value = "A012209532N"
$wsData.Range("C3").Value = "Paul Chotomolo Mirikwa"
$wsData.Range("D3").Value = "NOTICE UNDER SECTION 29 OF THE TAX PROCEDURES ACT, 2015"
$wsData.Range("E3").Value = "74,468.80"
$wsData.Range("H3").Value = "Franciscar Nyangweta"
$wsData.Range("I3").Value = "NAITIRI"
$wsData.Range("J3").Value = "2025-09-26 11:22:11"
$wsData.Range("L3").Value = 5

# Row 4
$wsData.Range("A4").Value = "04th September, 2025"
$wsData.Range("B4").Value = "A004578892U"
$wsData.Range("C4").Value = "JESSY KAGONDU WAMBUGU"
$wsData.Range("D4").Value = "NOTICE UNDER SECTION 29 OF THE TAX PROCEDURES ACT, 2015"
$wsData.Range("E4").Value = "118,561.81"
$wsData.Range("H4").Value = "Franciscar Nyangweta"
$wsData.Range("I4").Value = "KITALE"
$wsData.Range("J4").Value = "2025-09-26 10:42:57"
$wsData.Range("L4").Value = 5

# Row 5
$wsData.Range("A5").Value = "04TH September, 2025"
$wsData.Range("B5").Value = "A018905312S"
$wsData.Range("C5").Value = "Daisy Jepkosgei Biwott"
$wsData.Range("E5").Value = "1,348,612.53"
$wsData.Range("J5").Value = "2025-09-26 10:37:58"
$wsData.Range("L5").Value = 1

# Row 6 (new)
$wsData.Range("A6").Value = "4th September, 2025"
$wsData.Range("B6").Value = "A009775891W"
$wsData.Range("C6").Value = "Ezekiel Kipserem Korir"
$wsData.Range("D6").Value = "NOTICE UNDER SECTION 29 OF THE TAX PROCEDURES ACT, 2015"
$wsData.Range("E6").Value = "238,640.79"
$wsData.Range("G6").Value = 2024
$wsData.Range("H6").Value = "Franciscar Nyangweta"
$wsData.Range("I6").Value = "KITALE"
$wsData.Range("J6").Value = "2025-09-26 10:37:58"
$wsData.Range("K6").Value = "multi_format_extractor"
$wsData.Range("L6").Value = 2

# Row 7 (new)
$wsData.Range("A7").Value = "04th September, 2025"
$wsData.Range("B7").Value = "A005977112Z"
$wsData.Range("C7").Value = "James Mutoro Kitui"
$wsData.Range("D7").Value = "NOTICE UNDER SECTION 29 OF THE TAX PROCEDURES ACT, 2015"
$wsData.Range("E7").Value = "68,547.16"
$wsData.Range("G7").Value = 2024
$wsData.Range("H7").Value = "Franciscar Nyangweta"
$wsData.Range("I7").Value = "KITALE"
$wsData.Range("J7").Value = "2025-09-26 10:37:58"
$wsData.Range("K7").Value = "multi_format_extractor"
$wsData.Range("L7").Value = 3

# Restore the default (Normal) style on the amount column now that the
# text has been committed, so no stray "Text" number format lingers on
# these cells.
$wsData.Range("E2:E7").Style = "Normal"

# --- Database_Summary sheet --------------------------------------------
$wsSummary.Range("B2").Value = 6
$wsSummary.Range("B3").Value = "2025-09-26 11:25:36"
$wsSummary.Range("B5").Value = 0
$wsSummary.Range("B6").Value = "04TH September, 2025"
$wsSummary.Range("B7").Value = "4th September, 2025"
$wsSummary.Range("B8").Value = 6
$wsSummary.Range("B9").Value = 3
